$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.908.54'
$ws.Range("E2").Value = '  -2.01%  '
$ws.Range("D3").Value = '1.653.26'
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Formula = '="311.65"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").Formula = '="1.001"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Formula = '="0.3899"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  -1.24%  '
$ws.Range("D8").Formula = '="0.3827"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  -2.59%  '
$ws.Range("D9").Formula = '="51.69"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("D10").Formula = '="1.346"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  -3.31%  '
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Formula = '="0.08468"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  -1.04%  '
$ws.Range("D13").Formula = '="23.98"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  -2.00%  '
$ws.Range("D14").Formula = '="7.044"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  -3.19%  '
$ws.Range("D15").Formula = '="8.043"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").Formula = '="0.00001316"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  -1.51%  '
$ws.Range("D17").Value = '1.656.96'
$ws.Range("E17").Value = '  -0.57%  '
$ws.Range("D18").Formula = '="94.24"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  -0.67%  '
$ws.Range("D19").Formula = '="0.06993"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  -0.67%  '
$ws.Range("D20").Formula = '="19.64"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  -4.58%  '
$ws.Range("D21").Formula = '="6.977"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").Formula = '="1.000"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").Value = '23.888.17'
$ws.Range("E24").Value = '  -2.14%  '
$ws.Range("D25").Formula = '="2.450"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  -1.75%  '
$ws.Range("D26").Formula = '="2.955"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  -3.83%  '
$ws.Range("E27").Value = '  -2.06%  '
$ws.Range("D28").Formula = '="153.26"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  -2.45%  '
$ws.Range("D29").Formula = '="5.428"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").Formula = '="137.62"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  -3.46%  '
$ws.Range("D31").Formula = '="7.747"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  -2.65%  '
$ws.Range("E32").Value = '  -2.67%  '
$ws.Range("D33").Value = '1.838.86'
$ws.Range("E33").Value = '  -0.54%  '
$ws.Range("D34").Formula = '="1.007"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  -5.33%  '
$ws.Range("D35").Formula = '="0.08148"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -1.29%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Formula = '="0.02910"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  -6.36%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Formula = '="6.672"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  -3.31%  '
$ws.Range("D38").Formula = '="10.75"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  -3.20%  '
$ws.Range("D39").Formula = '="0.2678"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  -3.03%  '
$ws.Range("D40").Formula = '="0.09119"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  -1.56%  '
$ws.Range("D41").Formula = '="13.58"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -0.81%  '
$ws.Range("D42").Formula = '="0.7570"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  -1.51%  '
$ws.Range("D43").Formula = '="1.422"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  -1.65%  '
$ws.Range("D44").Formula = '="16.47"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  -0.43%  '
$ws.Range("D45").Formula = '="0.6937"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  -2.16%  '
$ws.Range("D46").Formula = '="2.450"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -3.63%  '
$ws.Range("D47").Formula = '="4.116"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("D48").Formula = '="0.9998"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("D49").Formula = '="0.08280"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  -1.85%  '
$ws.Range("D50").Formula = '="133.70"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  -2.24%  '
$ws.Range("D51").Formula = '="1.233"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  -2.64%  '
$excel.CutCopyMode = 0
